$d = $word.ActiveDocument

# The document has two TODO bullet paragraphs that must be removed entirely:
#   "Remove MetricDefinition specialization for Memory, Load, etc."
#   "Uncomment // @Test and make all test pass"
# The trailing "_GoBack" bookmark (which lives at the end of the second
# bullet, right before its own paragraph mark) must survive and end up on
# the empty paragraph that used to precede the two bullets.

# Locate the first bullet paragraph and the plain paragraph right before it.
$prev = $null
$target1 = $null
$prevCandidate = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Remove MetricDefinition*") {
        $target1 = $p
        $prev = $prevCandidate
        break
    }
    $prevCandidate = $p
}

# Overwrite the first bullet's content with just the plain (non-list)
# paragraph properties plus the bookmark - this is what should remain once
# both bullets are gone.
$bookmarkXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target1.Range.InsertXML($bookmarkXml)

# The following paragraph is now the second bullet ("Uncomment // @Test...");
# delete it completely (its text and its own paragraph mark).
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Uncomment // @Test*") {
        $target2 = $p
        break
    }
}
$target2.Range.Delete()

# Merge the placeholder paragraph (now holding just the Garamond rPr and the
# _GoBack bookmark) into the previously-empty paragraph ahead of it, by
# deleting that preceding paragraph's own mark. Only one paragraph remains
# at that spot afterward, carrying the bookmark.
$mark = $d.Range($prev.Range.End - 1, $prev.Range.End)
$mark.Delete()
